$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1492829.6
$ws.Range("I9").Value = 1865937.8
$ws.Range("K9").Value = 1865937.8
$ws.Range("M9").Value = -1865768.8
$ws.Range("H15").Value = 2699.375
$ws.Range("I15").Value = 2699.375
$ws.Range("K15").Value = 8098.125
$ws.Range("M15").Value = -7929.125
$ws.Range("H33").Value = 82.111115
$ws.Range("I33").Value = 82.111115
$ws.Range("K33").Value = 82.111115
$ws.Range("M33").Value = 146.888885
$ws.Range("H39").Value = 424.85715
$ws.Range("I39").Value = 325.33334
$ws.Range("K39").Value = 976.0000200000001
$ws.Range("M39").Value = -680.0000200000001
$ws.Range("H40").Value = 4471.857
$ws.Range("I40").Value = 3874.75
$ws.Range("J40").Value = 5268
$ws.Range("K40").Value = 3874.75
$ws.Range("L40").Value = 5268
$ws.Range("M40").Value = -3699.75
$ws.Range("N40").Value = -5618
$ws.Range("H86").Value = 8998.6
$ws.Range("I86").Value = 8748.25
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 8748.25
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -7625.25
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 8998.6
$ws.Range("I89").Value = 8748.25
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 43741.25
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -38125.25
$ws.Range("N89").Value = -61232
$ws.Range("H94").Value = 761.75
$ws.Range("I94").Value = 761.75
$ws.Range("K94").Value = 761.75
$ws.Range("M94").Value = -310.75
$ws.Range("H107").Value = 603.8125
$ws.Range("I107").Value = 539.7143
$ws.Range("J107").Value = 1052.5
$ws.Range("K107").Value = 539.7143
$ws.Range("L107").Value = 1052.5
$ws.Range("M107").Value = 1380.2857
$ws.Range("N107").Value = -4892.5
$ws.Range("H135").Value = 2573.125
$ws.Range("I135").Value = 1764.6666
$ws.Range("K135").Value = 15881.9994
$ws.Range("M135").Value = -13346.9994
$ws.Range("H137").Value = 1472232.8
$ws.Range("I137").Value = 2001705.5
$ws.Range("J137").Value = 1474.8889
$ws.Range("K137").Value = 6005116.5
$ws.Range("L137").Value = 4424.6667
$ws.Range("M137").Value = -6002566.5
$ws.Range("N137").Value = -9524.6667
$ws.Range("H138").Value = 2124.5881
$ws.Range("I138").Value = 583.4545000000001
$ws.Range("J138").Value = 4950
$ws.Range("K138").Value = 1750.3635
$ws.Range("L138").Value = 14850
$ws.Range("M138").Value = 3389.6365
$ws.Range("N138").Value = -25130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3530587.5
$ws.Range("I32").Value = 1687002.5
$ws.Range("K32").Value = 1687002.5
$ws.Range("M32").Value = -1686715.5
$ws.Range("H74").Value = 1699.44
$ws.Range("I74").Value = 1127.8823
$ws.Range("K74").Value = 1127.8823
$ws.Range("M74").Value = -253.8823
$ws.Range("H77").Value = 1699.44
$ws.Range("I77").Value = 1127.8823
$ws.Range("K77").Value = 5639.4115
$ws.Range("M77").Value = -1271.4115
$ws.Range("H122").Value = 3619.7856
$ws.Range("I122").Value = 3881.182
$ws.Range("J122").Value = 2661.3333
$ws.Range("K122").Value = 11643.546
$ws.Range("L122").Value = 7983.999899999999
$ws.Range("M122").Value = -9193.545999999998
$ws.Range("N122").Value = -12883.9999
$ws.Range("H132").Value = 1903.4
$ws.Range("I132").Value = 950.1111
$ws.Range("K132").Value = 2850.3333
$ws.Range("M132").Value = -320.3332999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H106").Value = 28932.125
$ws.Range("J106").Value = 28932.125
$ws.Range("L106").Value = 28932.125
$ws.Range("N106").Value = -31456.125
$ws.Range("H134").Value = 3233.9167
$ws.Range("I134").Value = 2521.6
$ws.Range("K134").Value = 7564.799999999999
$ws.Range("M134").Value = -5029.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1946.8889
$ws.Range("J58").Value = 2160.818
$ws.Range("L58").Value = 2160.818
$ws.Range("N58").Value = -2566.818
$ws.Range("H107").Value = 3847230.2
$ws.Range("I107").Value = 7143601.5
$ws.Range("J107").Value = 1464.1666
$ws.Range("K107").Value = 7143601.5
$ws.Range("L107").Value = 1464.1666
$ws.Range("M107").Value = -7141681.5
$ws.Range("N107").Value = -5304.1666
$ws.Range("H122").Value = 3012.8125
$ws.Range("J122").Value = 4104
$ws.Range("L122").Value = 12312
$ws.Range("N122").Value = -17212
$ws.Range("H134").Value = 4336.7188
$ws.Range("I134").Value = 4679.2
$ws.Range("J134").Value = 3113.5715
$ws.Range("K134").Value = 14037.6
$ws.Range("L134").Value = 9340.7145
$ws.Range("M134").Value = -11502.6
$ws.Range("N134").Value = -14410.7145
$ws.Range("H136").Value = 1946.8889
$ws.Range("J136").Value = 2160.818
$ws.Range("L136").Value = 6482.454000000001
$ws.Range("N136").Value = -11582.454
$ws.Range("H138").Value = 69165.5
$ws.Range("J138").Value = 69165.5
$ws.Range("L138").Value = 69165.5
$ws.Range("N138").Value = -79445.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5424603
$ws.Range("I4").Value = 5693775.5
$ws.Range("K4").Value = 17081326.5
$ws.Range("M4").Value = -17081214.5
$ws.Range("H113").Value = 1213.1111
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 1239.875
$ws.Range("K113").Value = 2997
$ws.Range("L113").Value = 3719.625
$ws.Range("M113").Value = -827
$ws.Range("N113").Value = -8059.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H107").Value = 1726.2307
$ws.Range("J107").Value = 2318.6667
$ws.Range("L107").Value = 2318.6667
$ws.Range("N107").Value = -6158.6667
$ws.Range("H113").Value = 1954.2106
$ws.Range("I113").Value = 1869.25
$ws.Range("K113").Value = 1869.25
$ws.Range("M113").Value = 300.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2305
$ws.Range("I12").Value = 2266.6667
$ws.Range("J12").Value = 2343.3333
$ws.Range("K12").Value = 2266.6667
$ws.Range("L12").Value = 2343.3333
$ws.Range("M12").Value = -2096.6667
$ws.Range("N12").Value = -2683.3333
$ws.Range("H22").Value = 1307.5714
$ws.Range("I22").Value = 1358.8334
$ws.Range("K22").Value = 1358.8334
$ws.Range("M22").Value = -1063.8334
$ws.Range("H27").Value = 1307.5714
$ws.Range("I27").Value = 1358.8334
$ws.Range("K27").Value = 1358.8334
$ws.Range("M27").Value = -1251.8334
$ws.Range("H46").Value = 2928.4285
$ws.Range("J46").Value = 7500
$ws.Range("L46").Value = 7500
$ws.Range("N46").Value = -7876
$ws.Range("H55").Value = 502.375
$ws.Range("I55").Value = 259.8889
$ws.Range("K55").Value = 259.8889
$ws.Range("M55").Value = -86.88889999999998
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 3964.8462
$ws.Range("J122").Value = 6147.6665
$ws.Range("L122").Value = 18442.9995
$ws.Range("N122").Value = -23342.9995
$ws.Range("H132").Value = 3622.125
$ws.Range("I132").Value = 3709.5334
$ws.Range("J132").Value = 3476.4443
$ws.Range("K132").Value = 11128.6002
$ws.Range("L132").Value = 10429.3329
$ws.Range("M132").Value = -8598.600199999999
$ws.Range("N132").Value = -15489.3329
$ws.Range("H136").Value = 5069.2383
$ws.Range("I136").Value = 5144.0586
$ws.Range("J136").Value = 4751.25
$ws.Range("K136").Value = 15432.1758
$ws.Range("L136").Value = 14253.75
$ws.Range("M136").Value = -12882.1758
$ws.Range("N136").Value = -19353.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H122").Value = 13890816
$ws.Range("I122").Value = 2011.75
$ws.Range("K122").Value = 6035.25
$ws.Range("M122").Value = -3585.25
$ws.Range("H126").Value = 13259.272
$ws.Range("J126").Value = 4666
$ws.Range("L126").Value = 13998
$ws.Range("N126").Value = -18938
$ws.Range("H132").Value = 2767.5334
$ws.Range("I132").Value = 2764.2083
$ws.Range("J132").Value = 2780.8333
$ws.Range("K132").Value = 8292.624899999999
$ws.Range("L132").Value = 8342.499899999999
$ws.Range("M132").Value = -5762.624899999999
$ws.Range("N132").Value = -13402.4999
$ws.Range("H135").Value = 289999.5
$ws.Range("J135").Value = 289999.5
$ws.Range("L135").Value = 289999.5
$ws.Range("N135").Value = -300139.5
$ws.Range("H136").Value = 3666.6667
$ws.Range("I136").Value = 3666.6667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11000.0001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8450.000100000001
$ws.Range("N136").ClearContents()
